# Fruta / hortaliza, semanal
# Insert a new data row for "Ají" (Inferno, Primera, Provincia de Limarí,
# $/caja 15 kilos) right before the existing row 308 (Fecha 2022-12-05),
# pushing that row and everything below it down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a whole new row at position 308 - this shifts old rows 308..385
# down to 309..386 and copies formatting (styles) from the row above.
$ws.Rows.Item(308).Insert()

# Populate the newly-inserted row 308 with the new record.
$ws.Cells.Item(308, 1).Value2  = 8
$ws.Cells.Item(308, 2).Value2  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(308, 3).Value2  = "Coquimbo"
$ws.Cells.Item(308, 4).Value2  = 44985
$ws.Cells.Item(308, 5).Value2  = 4
$ws.Cells.Item(308, 6).Value2  = 100112021
$ws.Cells.Item(308, 7).Value2  = "Ají"
$ws.Cells.Item(308, 8).Value2  = "Inferno"
$ws.Cells.Item(308, 9).Value2  = "Primera"
$ws.Cells.Item(308, 10).Value2 = 460
$ws.Cells.Item(308, 11).Value2 = 10500
$ws.Cells.Item(308, 12).Value2 = 11000
$ws.Cells.Item(308, 13).Value2 = 10750
$ws.Cells.Item(308, 14).Value2 = "$/caja 15 kilos"
$ws.Cells.Item(308, 15).Value2 = "Provincia de Limarí"
$ws.Cells.Item(308, 16).Value2 = 717
$ws.Cells.Item(308, 17).Value2 = 15
$ws.Cells.Item(308, 18).Value2 = "Hortaliza"
